$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1910.0227
$ws.Range("I132").Value = 1818.1842
$ws.Range("K132").Value = 5454.5526
$ws.Range("M132").Value = -2924.5526

$ws.Range("H135").Value = 714962.8
$ws.Range("I135").Value = 714962.8
$ws.Range("K135").Value = 6434665.2
$ws.Range("M135").Value = -6432130.2

$ws.Range("H137").Value = 3973.641
$ws.Range("I137").Value = 4981.091
$ws.Range("J137").Value = 3577.8572
$ws.Range("K137").Value = 14943.273
$ws.Range("L137").Value = 10733.5716
$ws.Range("M137").Value = -12393.273
$ws.Range("N137").Value = -15833.5716

$ws.Range("H141").Value = 1824.5151
$ws.Range("I141").Value = 1640.4333
$ws.Range("K141").Value = 4921.2999
$ws.Range("M141").Value = 258.7001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6571.7144
$ws.Range("I2").Value = 2671.3333
$ws.Range("K2").Value = 2671.3333
$ws.Range("M2").Value = -2558.3333

$ws.Range("H61").Value = 6280.4
$ws.Range("I61").Value = 4349.375
$ws.Range("J61").Value = 8487.286
$ws.Range("K61").Value = 4349.375
$ws.Range("L61").Value = 8487.286
$ws.Range("M61").Value = -4137.375
$ws.Range("N61").Value = -8911.286

$ws.Range("H116").Value = 6571.7144
$ws.Range("I116").Value = 2671.3333
$ws.Range("K116").Value = 2671.3333
$ws.Range("M116").Value = -377.3332999999998

$ws.Range("H132").Value = 3628.9487
$ws.Range("I132").Value = 1677.7241
$ws.Range("J132").Value = 9287.5
$ws.Range("K132").Value = 5033.1723
$ws.Range("L132").Value = 27862.5
$ws.Range("M132").Value = -2503.1723
$ws.Range("N132").Value = -32922.5

$ws.Range("H135").Value = 44429
$ws.Range("J135").Value = 44429
$ws.Range("L135").Value = 44429
$ws.Range("N135").Value = -54569

$ws.Range("H136").Value = 6280.4
$ws.Range("I136").Value = 4349.375
$ws.Range("J136").Value = 8487.286
$ws.Range("K136").Value = 13048.125
$ws.Range("L136").Value = 25461.858
$ws.Range("M136").Value = -10498.125
$ws.Range("N136").Value = -30561.858

$ws.Range("H139").Value = 70678.5
$ws.Range("J139").Value = 70678.5
$ws.Range("L139").Value = 70678.5
$ws.Range("N139").Value = -80958.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6571.7144
$ws.Range("I3").Value = 2671.3333
$ws.Range("K3").Value = 2671.3333
$ws.Range("M3").Value = -2557.3333

$ws.Range("H86").Value = 8656781
$ws.Range("I86").Value = 13211633
$ws.Range("J86").Value = 2564.4
$ws.Range("K86").Value = 13211633
$ws.Range("L86").Value = 2564.4
$ws.Range("M86").Value = -13210510
$ws.Range("N86").Value = -4810.4

$ws.Range("H89").Value = 8656781
$ws.Range("I89").Value = 13211633
$ws.Range("J89").Value = 2564.4
$ws.Range("K89").Value = 66058165
$ws.Range("L89").Value = 12822
$ws.Range("M89").Value = -66052549
$ws.Range("N89").Value = -24054

$ws.Range("H94").Value = 2106
$ws.Range("I94").Value = 1960.2
$ws.Range("J94").Value = 2178.9
$ws.Range("K94").Value = 1960.2
$ws.Range("L94").Value = 2178.9
$ws.Range("M94").Value = -1509.2
$ws.Range("N94").Value = -3080.9

$ws.Range("H134").Value = 5954.24
$ws.Range("I134").Value = 2465.1072
$ws.Range("K134").Value = 7395.321599999999
$ws.Range("M134").Value = -4860.321599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8952.604499999999
$ws.Range("I31").Value = 5293.25
$ws.Range("J31").Value = 13574.947
$ws.Range("K31").Value = 5293.25
$ws.Range("L31").Value = 13574.947
$ws.Range("M31").Value = -4998.25
$ws.Range("N31").Value = -14164.947

$ws.Range("H34").Value = 8952.604499999999
$ws.Range("I34").Value = 5293.25
$ws.Range("J34").Value = 13574.947
$ws.Range("K34").Value = 5293.25
$ws.Range("L34").Value = 13574.947
$ws.Range("M34").Value = -5091.25
$ws.Range("N34").Value = -13978.947

$ws.Range("H58").Value = 9863.440000000001
$ws.Range("J58").Value = 11367.45
$ws.Range("L58").Value = 11367.45
$ws.Range("N58").Value = -11773.45

$ws.Range("H60").Value = 38853.57
$ws.Range("I60").Value = 22000
$ws.Range("K60").Value = 22000
$ws.Range("M60").Value = -21489

$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws.Range("H99").Value = 4856.303
$ws.Range("I99").Value = 3587.5833
$ws.Range("K99").Value = 3587.5833
$ws.Range("M99").Value = -2089.5833

$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws.Range("H107").Value = 1389.3667
$ws.Range("I107").Value = 516.9091
$ws.Range("J107").Value = 1894.4736
$ws.Range("K107").Value = 516.9091
$ws.Range("L107").Value = 1894.4736
$ws.Range("M107").Value = 1403.0909
$ws.Range("N107").Value = -5734.4736

$ws.Range("H126").Value = 4856.303
$ws.Range("I126").Value = 3587.5833
$ws.Range("K126").Value = 10762.7499
$ws.Range("M126").Value = -8292.749899999999

$ws.Range("H132").Value = 5751.804
$ws.Range("I132").Value = 3914.0417
$ws.Range("K132").Value = 11742.1251
$ws.Range("M132").Value = -9212.125100000001

$ws.Range("H134").Value = 4458.4204
$ws.Range("J134").Value = 7854
$ws.Range("L134").Value = 23562
$ws.Range("N134").Value = -28632

$ws.Range("H136").Value = 9863.440000000001
$ws.Range("J136").Value = 11367.45
$ws.Range("L136").Value = 34102.35000000001
$ws.Range("N136").Value = -39202.35000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 54388.457
$ws.Range("I2").Value = 12556.5
$ws.Range("J2").Value = 100023.32
$ws.Range("K2").Value = 75339
$ws.Range("L2").Value = 600139.92
$ws.Range("M2").Value = -75226
$ws.Range("N2").Value = -600365.92

$ws.Range("H38").Value = 35.11111
$ws.Range("J38").Value = 22.666666
$ws.Range("L38").Value = 67.99999800000001
$ws.Range("N38").Value = -761.999998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 37638
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()

$ws.Range("H50").Value = 37638
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()

$ws.Range("H132").Value = 8650
$ws.Range("I132").Value = 1840
$ws.Range("K132").Value = 5520
$ws.Range("M132").Value = -2990

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2497.8333
$ws.Range("I22").Value = 1442
$ws.Range("J22").Value = 7777
$ws.Range("K22").Value = 1442
$ws.Range("L22").Value = 7777
$ws.Range("M22").Value = -1147
$ws.Range("N22").Value = -8367

$ws.Range("H27").Value = 2497.8333
$ws.Range("I27").Value = 1442
$ws.Range("J27").Value = 7777
$ws.Range("K27").Value = 1442
$ws.Range("L27").Value = 7777
$ws.Range("M27").Value = -1335
$ws.Range("N27").Value = -7991

$ws.Range("H68").Value = 5407.2383
$ws.Range("I68").Value = 3766
$ws.Range("K68").Value = 3766
$ws.Range("M68").Value = -3017

$ws.Range("H71").Value = 5407.2383
$ws.Range("I71").Value = 3766
$ws.Range("K71").Value = 18830
$ws.Range("M71").Value = -15086

$ws.Range("H82").Value = 644031.5
$ws.Range("I82").Value = 1283591.5
$ws.Range("J82").Value = 4471.4546
$ws.Range("K82").Value = 1283591.5
$ws.Range("L82").Value = 4471.4546
$ws.Range("M82").Value = -1283230.5
$ws.Range("N82").Value = -5193.4546

$ws.Range("H85").Value = 644031.5
$ws.Range("I85").Value = 1283591.5
$ws.Range("J85").Value = 4471.4546
$ws.Range("K85").Value = 1283591.5
$ws.Range("L85").Value = 4471.4546
$ws.Range("M85").Value = -1282343.5
$ws.Range("N85").Value = -6967.4546

$ws.Range("H93").Value = 5374.125
$ws.Range("I93").Value = 3828.9
$ws.Range("J93").Value = 7949.5
$ws.Range("K93").Value = 3828.9
$ws.Range("L93").Value = 7949.5
$ws.Range("M93").Value = -2580.9
$ws.Range("N93").Value = -10445.5

$ws.Range("H122").Value = 4747
$ws.Range("I122").Value = 2927.9333
$ws.Range("K122").Value = 8783.7999
$ws.Range("M122").Value = -6333.7999

$ws.Range("H132").Value = 9439805
$ws.Range("I132").Value = 16669422
$ws.Range("K132").Value = 50008266
$ws.Range("M132").Value = -50005736

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8248.25
$ws.Range("I62").Value = 8331
$ws.Range("K62").Value = 8331
$ws.Range("M62").Value = -7707

$ws.Range("H65").Value = 8248.25
$ws.Range("I65").Value = 8331
$ws.Range("K65").Value = 41655
$ws.Range("M65").Value = -38535

$ws.Range("H107").Value = 10753390
$ws.Range("I107").Value = 448.52942
$ws.Range("J107").Value = 23810534
$ws.Range("K107").Value = 1345.58826
$ws.Range("L107").Value = 71431602
$ws.Range("M107").Value = 574.41174
$ws.Range("N107").Value = -71435442
